# Updates countries & provincias Spain sheet ("Pais") to the new COVID snapshot.
# - refreshes Casos totales/Nuevos casos/Casos activos/Recuperados/Casos criticos/Muertes hoy/Muertes
#   for every country whose figures moved in the new pull,
# - re-sorts the handful of countries whose total-case rank crossed a neighbour
#   (Israel/Nigeria, Azerbaiyan/Guatemala, Argelia/Nepal, Guayana Francesa/Hungria, Laos/Santa Lucia),
# - bumps the "Datos actualizados" timestamp string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Timestamp footer/title cell
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 1 de Julio de 2020 a las 19:03'

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 2750688
$ws.Cells.Item(4, 3).Value = 22835
$ws.Cells.Item(4, 4).Value = 1146163
$ws.Cells.Item(4, 5).Value = 1474129
$ws.Cells.Item(4, 7).Value = 274
$ws.Cells.Item(4, 8).Value = 130396

# Row 5: Brasil
$ws.Cells.Item(5, 2).Value = 1426913
$ws.Cells.Item(5, 3).Value = 18428
$ws.Cells.Item(5, 5).Value = 576679
$ws.Cells.Item(5, 7).Value = 538
$ws.Cells.Item(5, 8).Value = 60194

# Row 7: India
$ws.Cells.Item(7, 2).Value = 604133
$ws.Cells.Item(7, 3).Value = 18341
$ws.Cells.Item(7, 4).Value = 359523
$ws.Cells.Item(7, 5).Value = 226783
$ws.Cells.Item(7, 7).Value = 417
$ws.Cells.Item(7, 8).Value = 17827

# Row 12: Italia
$ws.Cells.Item(12, 2).Value = 240760
$ws.Cells.Item(12, 3).Value = 182
$ws.Cells.Item(12, 4).Value = 190717
$ws.Cells.Item(12, 5).Value = 15255
$ws.Cells.Item(12, 7).Value = 21
$ws.Cells.Item(12, 8).Value = 34788

# Row 22: Canada
$ws.Cells.Item(22, 2).Value = 104271
$ws.Cells.Item(22, 3).Value = 67
$ws.Cells.Item(22, 4).Value = 67743
$ws.Cells.Item(22, 5).Value = 27913
$ws.Cells.Item(22, 7).Value = 24
$ws.Cells.Item(22, 8).Value = 8615

# Row 51: Israel
$ws.Cells.Item(51, 1).Value = 'Israel'
$ws.Cells.Item(51, 2).Value = 26021
$ws.Cells.Item(51, 3).Value = 777
$ws.Cells.Item(51, 4).Value = 17429
$ws.Cells.Item(51, 5).Value = 8271
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 321

# Row 52: Nigeria
$ws.Cells.Item(52, 1).Value = 'Nigeria'
$ws.Cells.Item(52, 2).Value = 25694
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 9746
$ws.Cells.Item(52, 5).Value = 15358
$ws.Cells.Item(52, 8).Value = 590

# Row 54: Kazajistan
$ws.Cells.Item(54, 4).Value = 13845
$ws.Cells.Item(54, 5).Value = 8275

# Row 57: Azerbaiyan
$ws.Cells.Item(57, 1).Value = 'Azerbaiyan'
$ws.Cells.Item(57, 2).Value = 18112
$ws.Cells.Item(57, 3).Value = 588
$ws.Cells.Item(57, 4).Value = 10061
$ws.Cells.Item(57, 5).Value = 7831
$ws.Cells.Item(57, 7).Value = 7
$ws.Cells.Item(57, 8).Value = 220

# Row 58: Guatemala
$ws.Cells.Item(58, 1).Value = 'Guatemala'
$ws.Cells.Item(58, 2).Value = 18096
$ws.Cells.Item(58, 3).Value = 687
$ws.Cells.Item(58, 4).Value = 3194
$ws.Cells.Item(58, 5).Value = 14129
$ws.Cells.Item(58, 7).Value = 27
$ws.Cells.Item(58, 8).Value = 773

# Row 59: Austria
$ws.Cells.Item(59, 1).Value = 'Austria'
$ws.Cells.Item(59, 2).Value = 17873
$ws.Cells.Item(59, 3).Value = 107
$ws.Cells.Item(59, 4).Value = 16491
$ws.Cells.Item(59, 5).Value = 677
$ws.Cells.Item(59, 8).Value = 705

# Row 60: Ghana
$ws.Cells.Item(60, 1).Value = 'Ghana'
$ws.Cells.Item(60, 2).Value = 17741
$ws.Cells.Item(60, 4).Value = 13268
$ws.Cells.Item(60, 5).Value = 4361
$ws.Cells.Item(60, 8).Value = 112

# Row 63: Argelia
$ws.Cells.Item(63, 1).Value = 'Argelia'
$ws.Cells.Item(63, 2).Value = 14272
$ws.Cells.Item(63, 3).Value = 365
$ws.Cells.Item(63, 4).Value = 10040
$ws.Cells.Item(63, 5).Value = 3312
$ws.Cells.Item(63, 7).Value = 8
$ws.Cells.Item(63, 8).Value = 920

# Row 64: Nepal
$ws.Cells.Item(64, 1).Value = 'Nepal'
$ws.Cells.Item(64, 2).Value = 14046
$ws.Cells.Item(64, 3).Value = 482
$ws.Cells.Item(64, 4).Value = 4656
$ws.Cells.Item(64, 5).Value = 9360
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = 30

# Row 69: Chequia
$ws.Cells.Item(69, 2).Value = 12006
$ws.Cells.Item(69, 3).Value = 52
$ws.Cells.Item(69, 4).Value = 7797
$ws.Cells.Item(69, 5).Value = 3860

# Row 94: Guayana Francesa
$ws.Cells.Item(94, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(94, 2).Value = 4268
$ws.Cells.Item(94, 3).Value = 264
$ws.Cells.Item(94, 4).Value = 1602
$ws.Cells.Item(94, 5).Value = 2650
$ws.Cells.Item(94, 8).Value = 16

# Row 95: Hungria
$ws.Cells.Item(95, 1).Value = 'Hungria'
$ws.Cells.Item(95, 2).Value = 4157
$ws.Cells.Item(95, 3).Value = 2
$ws.Cells.Item(95, 4).Value = 2714
$ws.Cells.Item(95, 5).Value = 857
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 8).Value = 586

# Row 130: Jordania
$ws.Cells.Item(130, 2).Value = 1133
$ws.Cells.Item(130, 3).Value = 1
$ws.Cells.Item(130, 4).Value = 886
$ws.Cells.Item(130, 5).Value = 238

# Row 135: Republica de Chipre
$ws.Cells.Item(135, 2).Value = 999
$ws.Cells.Item(135, 3).Value = 1
$ws.Cells.Item(135, 5).Value = 147

# Row 203: Laos
$ws.Cells.Item(203, 1).Value = 'Laos'

# Row 204: Santa Lucia
$ws.Cells.Item(204, 1).Value = 'Santa Lucia'
